$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scheduled data refresh: updated crypto prices / 1h volume deltas,
# and Hedera/Kaspa swapped ranking positions (rows 28-29).

# Force column D (Price) to retain text storage (values like "51.525.25" or
# "0.998" must stay text, matching how the sheet already stores them) by
# quoting the range as Text before writing, then restoring the default style
# so no extra number-format/style is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '51.525.25'
$ws.Cells.Item(2, 5).Value = '  +0.50%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.981.48'
$ws.Cells.Item(3, 5).Value = '  +2.11%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '0.998'
$ws.Cells.Item(4, 5).Value = '  -0.13%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '378.85'
$ws.Cells.Item(5, 5).Value = '  +2.44%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '104.83'
$ws.Cells.Item(6, 5).Value = '  +0.63%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.541'
$ws.Cells.Item(7, 5).Value = '  +0.06%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.06%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.595'
$ws.Cells.Item(9, 5).Value = '  +1.35%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '37.24'
$ws.Cells.Item(10, 5).Value = '  +1.34%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.02%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '0.0846'
$ws.Cells.Item(12, 5).Value = '  +1.46%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '3.447.03'
$ws.Cells.Item(13, 5).Value = '  +1.94%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '18.39'
$ws.Cells.Item(14, 5).Value = '  +0.12%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '7.60'
$ws.Cells.Item(15, 5).Value = '  +2.60%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.985.92'
$ws.Cells.Item(16, 5).Value = '  +2.27%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '0.967'
$ws.Cells.Item(17, 5).Value = '  +3.29%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '51.444.99'
$ws.Cells.Item(18, 5).Value = '  +0.45%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +2.74%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '7.43'
$ws.Cells.Item(20, 5).Value = '  +3.05%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '12.95'
$ws.Cells.Item(21, 5).Value = '  -0.04%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.0₃0963'
$ws.Cells.Item(22, 5).Value = '  +2.07%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '69.39'
$ws.Cells.Item(23, 5).Value = '  +1.43%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '262.22'
$ws.Cells.Item(24, 5).Value = '  +0.89%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '2.84'
$ws.Cells.Item(25, 5).Value = '  +5.65%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '8.24'
$ws.Cells.Item(26, 5).Value = '  +16.47%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '7.67'
$ws.Cells.Item(27, 5).Value = '  +22.84%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Hedera'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(28, 4).Value = '0.115'
$ws.Cells.Item(28, 5).Value = '  +11.68%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Kaspa'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(29, 4).Value = '0.170'
$ws.Cells.Item(29, 5).Value = '  -2.59%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -0.07%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '25.91'
$ws.Cells.Item(31, 5).Value = '  +0.52%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '9.88'
$ws.Cells.Item(32, 5).Value = '  -0.28%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '35.17'
$ws.Cells.Item(33, 5).Value = '  +1.02%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -2.07%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '51.07'
$ws.Cells.Item(35, 5).Value = '  +0.43%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '0.0446'
$ws.Cells.Item(36, 5).Value = '  +4.82%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.03%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '3.03'
$ws.Cells.Item(38, 5).Value = '  +0.46%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '17.22'
$ws.Cells.Item(39, 5).Value = '  +0.57%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '2.59'
$ws.Cells.Item(40, 5).Value = '  -2.98%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.34%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +2.26%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '125.58'
$ws.Cells.Item(43, 5).Value = '  +5.59%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '21.73'
$ws.Cells.Item(44, 5).Value = '  -2.06%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '0.283'
$ws.Cells.Item(45, 5).Value = '  +17.27%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '2.04'
$ws.Cells.Item(46, 5).Value = '  -1.51%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '2.39'
$ws.Cells.Item(47, 5).Value = '  +3.52%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '2.035.42'
$ws.Cells.Item(48, 5).Value = '  +0.73%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '3.24'
$ws.Cells.Item(49, 5).Value = '  +1.96%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '0.0336'
$ws.Cells.Item(50, 5).Value = '  +8.13%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '58.40'
$ws.Cells.Item(51, 5).Value = '  +2.67%  '

# Restore the default cell style on column D so the underlying style index
# is unchanged (only the stored value/type differs from before).
$ws.Range("D2:D51").Style = "Normal"
